$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Member table (rows 10-13): add student names + student id ---
$ws.Range("C10").Value = "Vương Thế Minh Thăng"
$ws.Range("D10").Value = "PS09070"
$ws.Range("C11").Value = "Võ Thành Long"
$ws.Range("C12").Value = "Hoàng Hồng Sơn"
$ws.Range("C13").Value = "Phạm Lê Huy"

# --- Task table (rows 16-22): durations, assignees, completion % ---
$ws.Range("C16").Value = 3
$ws.Range("G16").Value = "SV1, SV2"
$ws.Range("H16").Value = 0.8

$ws.Range("C17").Value = 2
$ws.Range("G17").Value = "SV1, SV2"
$ws.Range("H17").Value = 0.7

$ws.Range("C18").Value = 5
$ws.Range("G18").Value = "SV1, SV2"
$ws.Range("H18").Value = 0.9

$ws.Range("C19").Value = 4
$ws.Range("G19").Value = "SV1, SV2"
$ws.Range("H19").Value = 0.6

$ws.Range("C20").Value = 2
$ws.Range("G20").Value = "SV1"
$ws.Range("H20").Value = 1

$ws.Range("C21").Value = 3
$ws.Range("G21").Value = "SV2"
$ws.Range("H21").Value = 0.6

$ws.Range("C22").Value = 1
$ws.Range("G22").Value = "SV1, SV2"

# Percent number format for the "Hoàn thành (%)" column on the filled rows
$ws.Range("H16:H21").NumberFormat = "0%"

# Move the selection to match where the author ended up editing
$ws.Range("E17").Select()
